{"js": "// The default (primary) page header currently reads \"S\u00e4ulen, Freude,\"\n// and needs to be updated to \"Regenbogen, Bild, Hand, S\u00e4ulen\".\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst section = sections.items[0];\nconst header = section.getHeader(Word.HeaderFooterType.primary);\n\nconst results = header.search(\"S\u00e4ulen, Freude,\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Regenbogen, Bild, Hand, S\u00e4ulen\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The default (primary) page header currently reads \"S\u00e4ulen, Freude,\"\n# and needs to be updated to \"Regenbogen, Bild, Hand, S\u00e4ulen\".\n$d = $word.ActiveDocument\n$sec = $d.Sections.Item(1)\n$hdr = $sec.Headers.Item(1)   # wdHeaderFooterPrimary\n$rng = $hdr.Range\n\n$rng.Find.ClearFormatting()\n$rng.Find.Execute(\"S\u00e4ulen, Freude,\", $false, $false, $false, $false, $false, $true, 1, $false, \"Regenbogen, Bild, Hand, S\u00e4ulen\", 2)\n"}
